# Applies the "Final Form of announcement" revision to table7_1:
# updated N / CR(%) / AR(%) figures for each PEIS bucket, plus the
# significance-annotated AR value that moves from row 4 to row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ALL / N)
$ws.Range("D2").Value = 1338
$ws.Range("E2").Value = 1956
$ws.Range("G2").Value = 2538
$ws.Range("H2").Value = 3241

# Row 3 (ALL / CR (%))
$ws.Range("D3").Value = 0.1522
$ws.Range("E3").Value = 0.1564
$ws.Range("G3").Value = 0.2087
$ws.Range("H3").Value = 0.183

# Row 4 (ALL / AR (%)) - I4 switches from text "Hold" to a numeric value
$ws.Range("C4").Value = 0.16
$ws.Range("D4").Value = 0.1528
$ws.Range("E4").Value = 0.1784
$ws.Range("G4").Value = 0.1905
$ws.Range("H4").Value = 0.1276
$ws.Range("I4").Value = -0.0252

# Row 5 (PEIS1 / N)
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6 (PEIS1 / CR (%))
$ws.Range("D6").Value = -0.3629
$ws.Range("E6").Value = -0.397
$ws.Range("F6").Value = 0.2077

# Row 7 (PEIS1 / AR (%))
$ws.Range("C7").Value = 0.33
$ws.Range("H7").Value = -0.3076
$ws.Range("I7").Value = -0.3076

# Row 8 (PEIS2 / N)
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 37
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 54

# Row 9 (PEIS2 / CR (%))
$ws.Range("C9").Value = 0.14
$ws.Range("D9").Value = 0.1904
$ws.Range("E9").Value = -0.0047
$ws.Range("F9").Value = 0.1772
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0.1445

# Row 10 (PEIS2 / AR (%))
$ws.Range("C10").Value = 0.1
$ws.Range("D10").Value = 0.2293
$ws.Range("E10").Value = 0.122
$ws.Range("F10").Value = 0.1919
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.0673
$ws.Range("I10").Value = -0.162

# Row 11 (PEIS3 / N)
$ws.Range("D11").Value = 324
$ws.Range("E11").Value = 355
$ws.Range("G11").Value = 460
$ws.Range("H11").Value = 621

# Row 12 (PEIS3 / CR (%))
$ws.Range("D12").Value = 0.0791
$ws.Range("E12").Value = 0.1234
$ws.Range("G12").Value = 0.0697
$ws.Range("H12").Value = 0.0637

# Row 13 (PEIS3 / AR (%)) - I13 switches from a numeric value to the
# significance-annotated text that used to live in I4 ("-0.0722**")
$ws.Range("C13").Value = 0.09
$ws.Range("D13").Value = 0.1294
$ws.Range("E13").Value = 0.097
$ws.Range("G13").Value = 0.1539
$ws.Range("H13").Value = 0.0572
$ws.Range("I13").Value = "-0.0722**"

# Row 14 (PEIS4 / N)
$ws.Range("D14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 9

# Row 15 (PEIS4 / CR (%))
$ws.Range("C15").Value = -0.04
$ws.Range("D15").Value = -0.0903
$ws.Range("E15").Value = -0.0299
$ws.Range("F15").Value = -0.3267
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0.0013

# Row 16 (PEIS4 / AR (%))
$ws.Range("C16").Value = 0.03
$ws.Range("D16").Value = -0.0195
$ws.Range("E16").Value = -0.515
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.3108
$ws.Range("I16").Value = 0.3303
